$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, new Price (D) text (or $null if unchanged),
# new Volume(1h) (E) text, and whether D must be forced to stay text
# (Excel would otherwise silently coerce plain-looking numbers like
# "1.00" or "596.14" into numeric values, losing the trailing/format digits).
$changes = @(
    @{Row=2; D="67.652.61"; E="  -1.12%  "; ForceText=$false}
    @{Row=3; D="3.772.12"; E="  -2.00%  "; ForceText=$false}
    @{Row=4; D="1.00"; E="  -0.04%  "; ForceText=$true}
    @{Row=5; D="596.14"; E="  -0.96%  "; ForceText=$true}
    @{Row=6; D="168.44"; E="  -0.37%  "; ForceText=$true}
    @{Row=7; D="3.771.84"; E="  -2.01%  "; ForceText=$false}
    @{Row=8; D=$null; E="  +0.01%  "; ForceText=$false}
    @{Row=9; D=$null; E="  -0.98%  "; ForceText=$false}
    @{Row=10; D=$null; E="  -0.94%  "; ForceText=$false}
    @{Row=11; D=$null; E="  +0.57%  "; ForceText=$false}
    @{Row=12; D=$null; E="  -1.33%  "; ForceText=$false}
    @{Row=13; D="0.0000277"; E="  +3.89%  "; ForceText=$true}
    @{Row=14; D="36.27"; E="  -2.24%  "; ForceText=$true}
    @{Row=15; D="4.408.12"; E="  -1.92%  "; ForceText=$false}
    @{Row=16; D="3.777.49"; E="  -1.80%  "; ForceText=$false}
    @{Row=17; D="18.51"; E="  -0.16%  "; ForceText=$true}
    @{Row=18; D="67.582.39"; E="  -1.35%  "; ForceText=$false}
    @{Row=19; D="7.17"; E="  -2.83%  "; ForceText=$true}
    @{Row=20; D=$null; E="  +0.68%  "; ForceText=$false}
    @{Row=21; D="10.47"; E="  -6.70%  "; ForceText=$true}
    @{Row=22; D="466.44"; E="  -0.93%  "; ForceText=$true}
    @{Row=23; D=$null; E="  -2.35%  "; ForceText=$false}
    @{Row=24; D=$null; E="  -7.86%  "; ForceText=$false}
    @{Row=25; D="83.66"; E="  +0.19%  "; ForceText=$true}
    @{Row=26; D="2.19"; E="  -1.97%  "; ForceText=$true}
    @{Row=27; D="12.10"; E="  -0.11%  "; ForceText=$true}
    @{Row=28; D="10.30"; E="  +0.48%  "; ForceText=$true}
    @{Row=29; D=$null; E="  -0.10%  "; ForceText=$false}
    @{Row=30; D=$null; E="  -1.74%  "; ForceText=$false}
    @{Row=31; D="3.926.12"; E="  -1.82%  "; ForceText=$false}
    @{Row=32; D="7.59"; E="  -1.74%  "; ForceText=$true}
    @{Row=33; D="30.47"; E="  -3.25%  "; ForceText=$true}
    @{Row=34; D=$null; E="  -3.91%  "; ForceText=$false}
    @{Row=35; D="9.10"; E="  -2.80%  "; ForceText=$true}
    @{Row=36; D="3.737.70"; E="  -2.03%  "; ForceText=$false}
    @{Row=37; D="3.69"; E="  -2.54%  "; ForceText=$true}
    @{Row=38; D=$null; E="  -1.34%  "; ForceText=$false}
    @{Row=39; D=$null; E="  -2.01%  "; ForceText=$false}
    @{Row=40; D="0.137"; E="  -1.79%  "; ForceText=$true}
    @{Row=41; D=$null; E="  -2.71%  "; ForceText=$false}
    @{Row=42; D="0.999"; E="  -0.05%  "; ForceText=$true}
    @{Row=43; D=$null; E="  -1.58%  "; ForceText=$false}
    @{Row=45; D=$null; E="  -0.86%  "; ForceText=$false}
    @{Row=46; D="1.93"; E="  -2.66%  "; ForceText=$true}
    @{Row=47; D="45.68"; E="  -2.77%  "; ForceText=$true}
    @{Row=48; D="395.18"; E="  -5.44%  "; ForceText=$true}
    @{Row=49; D=$null; E="  -8.72%  "; ForceText=$false}
    @{Row=50; D="140.48"; E="  -1.11%  "; ForceText=$true}
    @{Row=51; D="39.18"; E="  +2.80%  "; ForceText=$true}
)

foreach ($change in $changes) {
    if ($null -ne $change.D) {
        $dCell = $ws.Range("D" + $change.Row)
        if ($change.ForceText) {
            # Pin the cell to text format so the numeric-looking string
            # ("1.00", "596.14", ...) round-trips verbatim instead of
            # being parsed into a Double, then drop the format override
            # again so no visible/applied formatting change remains.
            $dCell.NumberFormat = "@"
            $dCell.Value = $change.D
            $dCell.ClearFormats()
        } else {
            $dCell.Value = $change.D
        }
    }
    $ws.Range("E" + $change.Row).Value = $change.E
}
